# Updated symbol list with refreshed coin prices, volumes, and rankings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Address, $Value)
    $cell = $ws.Range($Address)
    # Force text storage so numeric-looking strings (prices, percentages)
    # keep their exact formatting (e.g. trailing zeros) instead of being
    # auto-coerced to numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

Set-TextValue "D2" '305.30'
Set-TextValue "E2" '2.62%'
Set-TextValue "D3" '44.27'
Set-TextValue "E3" '7.37%'
Set-TextValue "D4" '5.123'
Set-TextValue "E4" '2.40%'
Set-TextValue "D5" '0.07838'
Set-TextValue "E5" '4.05%'
Set-TextValue "B6" 'GateToken'
Set-TextValue "C6" 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D6" '4.432'
Set-TextValue "E6" '1.62%'
Set-TextValue "B7" 'FTXToken'
Set-TextValue "C7" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D7" '1.616'
Set-TextValue "E7" '2.99%'
Set-TextValue "B8" 'MXToken'
Set-TextValue "C8" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D8" '1.057'
Set-TextValue "E8" '13.69%'
Set-TextValue "B9" 'LiechtensteinCryptoassetsExchange'
Set-TextValue "C9" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D9" '0.1302'
Set-TextValue "E9" '7.18%'
Set-TextValue "B10" 'WazirX'
Set-TextValue "C10" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D10" '0.1864'
Set-TextValue "E10" '2.05%'
Set-TextValue "B11" 'MandalaExchangeToken'
Set-TextValue "C11" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D11" '0.09163'
Set-TextValue "E11" '3.46%'
Set-TextValue "B12" 'BitrueCoin'
Set-TextValue "C12" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D12" '0.04145'
Set-TextValue "E12" '1.84%'
Set-TextValue "B13" 'BitMartToken'
Set-TextValue "C13" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D13" '0.1045'
Set-TextValue "E13" '-0.69%'
Set-TextValue "B14" 'BitForexToken'
Set-TextValue "C14" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D14" '0.001298'
Set-TextValue "E14" '1.41%'
Set-TextValue "B15" 'TigerCash'
Set-TextValue "C15" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D15" '0.005875'
Set-TextValue "E15" '-0.18%'
Set-TextValue "B16" 'UpBots'
Set-TextValue "C16" 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue "D16" '0.007445'
Set-TextValue "E16" '1,899.16%'
Set-TextValue "B17" 'LEO'
Set-TextValue "C17" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.369'
Set-TextValue "E17" '0.69%'
Set-TextValue "D18" '2.344'
Set-TextValue "E18" '-2.66%'
Set-TextValue "D19" '0.3367'
Set-TextValue "E19" '2.24%'
Set-TextValue "D20" '8.023'
Set-TextValue "E20" '0.89%'
Set-TextValue "D21" '0.1372'
Set-TextValue "E21" '-3.18%'
Set-TextValue "E22" '-5.14%'
Set-TextValue "D23" '0.04176'
Set-TextValue "E23" '3.07%'
Set-TextValue "D24" '0.001274'
Set-TextValue "E24" '0.91%'
Set-TextValue "D25" '0.004453'
Set-TextValue "E25" '14.04%'
Set-TextValue "D26" '0.0001340'
Set-TextValue "E26" '9.14%'
Set-TextValue "D38" '0.02531'
Set-TextValue "E38" '4.44%'
Set-TextValue "D39" '0.05338'
Set-TextValue "E39" '2.42%'
Set-TextValue "D40" '0.005453'
Set-TextValue "E40" '-7.57%'
Set-TextValue "D41" '0.007794'
Set-TextValue "E41" '0.17%'
Set-TextValue "E42" '2.94%'
Set-TextValue "D43" '0.007330'
Set-TextValue "E43" '-0.40%'
Set-TextValue "D44" '0.008327'
Set-TextValue "E44" '6.35%'
Set-TextValue "D45" '0.3027'
Set-TextValue "E45" '1.76%'
Set-TextValue "D46" '0.00006674'
Set-TextValue "E46" '5.69%'
Set-TextValue "D47" '0.00000000745'
Set-TextValue "E48" '31.35%'
Set-TextValue "D49" '0.003976'
Set-TextValue "E49" '-5.20%'
Set-TextValue "D50" '0.00002087'
Set-TextValue "D51" '0.0001988'
